$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A2:D18")
$keyRng = $ws.Range("A2:A18")

$rng.Sort($keyRng, 1)
